$wb = $excel.ActiveWorkbook

# Update the unit value in yeni_degiskenler!B2 (72 -> 62); downstream formulas
# in yeni_otv!G2:G7 reference this cell and will recalc automatically.
$ws = $wb.Worksheets.Item("yeni_degiskenler")
$ws.Range("B2").Value = 62

# Move the active selection on this sheet from B6 to B3
$ws.Activate()
$ws.Range("B3").Select()
